$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells G1:N1
$headers = @("min", "q02", "q25", "median", "q75", "q98", "max", "avg")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = 7 + $i   # G = 7 ... N = 14
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $headers[$i]
}

# Copy the header style/format from F1 (bold, bordered, centered) to G1:N1
$ws.Range("F1").Copy()
$ws.Range("G1:N1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$data = New-Object 'object[,]' 27,8
$data[0,0] = -1
$data[0,1] = -1
$data[0,2] = 4
$data[0,3] = 7
$data[0,4] = 10
$data[0,5] = 21
$data[0,6] = 44
$data[0,7] = 7.829073905944824
$data[1,0] = -1
$data[1,1] = -1
$data[1,2] = 4
$data[1,3] = 7
$data[1,4] = 11
$data[1,5] = 23
$data[1,6] = 36
$data[1,7] = 7.932900428771973
$data[2,0] = -1
$data[2,1] = -1
$data[2,2] = 4
$data[2,3] = 7
$data[2,4] = 11
$data[2,5] = 22
$data[2,6] = 50
$data[2,7] = 7.747605800628662
$data[3,0] = 0
$data[3,1] = 0
$data[3,2] = 4
$data[3,3] = 6
$data[3,4] = 10
$data[3,5] = 22
$data[3,6] = 49
$data[3,7] = 7.5487380027771
$data[4,0] = 0
$data[4,1] = 0
$data[4,2] = 4
$data[4,3] = 7
$data[4,4] = 10
$data[4,5] = 22
$data[4,6] = 80
$data[4,7] = 7.633392810821533
$data[5,0] = 0
$data[5,1] = 0
$data[5,2] = 5
$data[5,3] = 8
$data[5,4] = 11
$data[5,5] = 23
$data[5,6] = 52
$data[5,7] = 8.692742347717285
$data[6,0] = -1
$data[6,1] = -1
$data[6,2] = 4
$data[6,3] = 6
$data[6,4] = 9
$data[6,5] = 20
$data[6,6] = 53
$data[6,7] = 6.805400371551514
$data[7,0] = -1
$data[7,1] = -1
$data[7,2] = 3
$data[7,3] = 6
$data[7,4] = 10
$data[7,5] = 30.22000000000025
$data[7,6] = 204
$data[7,7] = 7.962154696132597
$data[8,0] = 0
$data[8,1] = 0
$data[8,2] = 4
$data[8,3] = 7
$data[8,4] = 10
$data[8,5] = 19
$data[8,6] = 125
$data[8,7] = 7.236739515652688
$data[9,0] = -1
$data[9,1] = -1
$data[9,2] = 4
$data[9,3] = 6
$data[9,4] = 9
$data[9,5] = 20
$data[9,6] = 360
$data[9,7] = 7.104579792256846
$data[10,0] = 1.3
$data[10,1] = 1.3
$data[10,2] = 3.6
$data[10,3] = 5.4
$data[10,4] = 8.699999999999999
$data[10,5] = 21.8
$data[10,6] = 137.8
$data[10,7] = 7.007795536791313
$data[11,0] = 0
$data[11,1] = 0
$data[11,2] = 4
$data[11,3] = 7
$data[11,4] = 10
$data[11,5] = 22
$data[11,6] = 54
$data[11,7] = 7.904960632324219
$data[12,0] = 1
$data[12,1] = 1
$data[12,2] = 5
$data[12,3] = 7
$data[12,4] = 11
$data[12,5] = 21
$data[12,6] = 78
$data[12,7] = 8.091594696044922
$data[13,0] = 1
$data[13,1] = 1
$data[13,2] = 6
$data[13,3] = 9
$data[13,4] = 13
$data[13,5] = 24
$data[13,6] = 69
$data[13,7] = 9.649613380432129
$data[14,0] = 0
$data[14,1] = 0
$data[14,2] = 5
$data[14,3] = 7
$data[14,4] = 11
$data[14,5] = 25
$data[14,6] = 53
$data[14,7] = 8.487636566162109
$data[15,0] = -1
$data[15,1] = -1
$data[15,2] = 3
$data[15,3] = 5
$data[15,4] = 8
$data[15,5] = 20
$data[15,6] = 69
$data[15,7] = 5.902006149291992
$data[16,0] = 0
$data[16,1] = 0
$data[16,2] = 4
$data[16,3] = 7
$data[16,4] = 10
$data[16,5] = 21
$data[16,6] = 52
$data[16,7] = 7.549196720123291
$data[17,0] = 2
$data[17,1] = 2
$data[17,2] = 5
$data[17,3] = 9
$data[17,4] = 14
$data[17,5] = 24.45999908447266
$data[17,6] = 40
$data[17,7] = 10.21666622161865
$data[18,0] = -1
$data[18,1] = -1
$data[18,2] = 3
$data[18,3] = 6
$data[18,4] = 10
$data[18,5] = 20
$data[18,6] = 62
$data[18,7] = 6.956928253173828
$data[19,0] = -1
$data[19,1] = -1
$data[19,2] = 4
$data[19,3] = 6
$data[19,4] = 10
$data[19,5] = 33
$data[19,6] = 244
$data[19,7] = 8.657696597880648
$data[20,0] = 0
$data[20,1] = 0
$data[20,2] = 3
$data[20,3] = 6
$data[20,4] = 9
$data[20,5] = 19
$data[20,6] = 85
$data[20,7] = 6.686489058039962
$data[21,0] = 1
$data[21,1] = 1
$data[21,2] = 4
$data[21,3] = 7
$data[21,4] = 10
$data[21,5] = 20
$data[21,6] = 131
$data[21,7] = 7.690385530699666
$data[22,0] = 2.8512
$data[22,1] = 2.8512
$data[22,2] = 5.03
$data[22,3] = 7.117083333333333
$data[22,4] = 10.47625
$data[22,5] = 26.44700000000003
$data[22,6] = 175.495
$data[22,7] = 8.931984381063245
$data[23,0] = 0.9615833379576604
$data[23,1] = 0.9615833379576604
$data[23,2] = 4.528125002980232
$data[23,3] = 9.709750016530354
$data[23,4] = 16.31174997488658
$data[23,5] = 43.18033254623413
$data[23,6] = 202.4250005086263
$data[23,7] = 12.6525999722694
$data[24,0] = 0.7065541703402995
$data[24,1] = 0.7065541703402995
$data[24,2] = 4.619281297922134
$data[24,3] = 9.706145902474722
$data[24,4] = 18.08900010585785
$data[24,5] = 38.98978563308715
$data[24,6] = 148.2515411376953
$data[24,7] = 12.66643034488212
$data[25,0] = 0
$data[25,1] = 0
$data[25,2] = 3.355833331743876
$data[25,3] = 7.675625006357829
$data[25,4] = 26.57737493515015
$data[25,5] = 47.48729972839356
$data[25,6] = 96.08124923706055
$data[25,7] = 14.88347807317683
$data[26,0] = 0
$data[26,1] = 0
$data[26,2] = 0.8354166708886623
$data[26,3] = 5.08216667175293
$data[26,4] = 15.78683334986369
$data[26,5] = 41.69837659200033
$data[26,6] = 221.996084416906
$data[26,7] = 10.12769062769757

$rng = $ws.Range("G2:N28")
$rng.Value = $data
